# Append a new attendance record for "Deep Javya" as row 7, matching the
# layout/formatting already used by the table (rows 2-6: Date, Name, Email,
# In Time, Out Time).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 7 with row 6's formatting (number format / fonts / borders /
# alignment) by copying the formatted range down, then overwrite the values.
$ws.Range("A6:E6").Copy($ws.Range("A7:E7"))

$ws.Range("A7").Value = 45408
$ws.Range("B7").Value = "Deep Javya"
$ws.Range("C7").Value = "staff2it@adit.ac.in"
$ws.Range("D7").Value = "08:59AM"
$ws.Range("E7").Value = "05:05PM"

# The rest of the table uses a 15pt row height.
$ws.Rows(7).RowHeight = 15

# Leave the selection where the author left it after entering the new row.
$ws.Range("D7").Select() | Out-Null
